$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 3; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 14; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 20; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 42; I = 'ba'; J = 'Appreciation' }
    @{ Row = 48; I = 'ba'; J = 'Appreciation' }
    @{ Row = 49; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 54; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 56; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 57; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 65; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 66; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 69; I = 'ba'; J = 'Appreciation' }
    @{ Row = 75; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 76; I = 'ba'; J = 'Appreciation' }
    @{ Row = 109; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 118; I = 'ba'; J = 'Appreciation' }
    @{ Row = 123; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 137; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 139; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 140; I = '%'; J = 'Uninterpretable' }
    @{ Row = 141; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 143; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 161; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 163; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 165; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 185; I = 'ba'; J = 'Appreciation' }
    @{ Row = 194; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 207; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 219; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 220; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 251; I = 'ba'; J = 'Appreciation' }
    @{ Row = 253; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 261; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 274; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 280; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 296; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 309; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 315; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 321; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 345; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 355; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 359; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 360; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 362; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 372; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 379; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 391; I = 'ba'; J = 'Appreciation' }
    @{ Row = 393; I = '%'; J = 'Uninterpretable' }
    @{ Row = 397; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 405; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 411; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 415; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 426; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 442; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 446; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 449; I = 'ba'; J = 'Appreciation' }
    @{ Row = 452; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 453; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 457; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 468; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 479; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 494; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 507; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 509; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 510; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 537; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 548; I = 'ba'; J = 'Appreciation' }
    @{ Row = 556; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"